# Fixed Stimulus Absolute Timestamps
# Rename each task-order sheet tab to reflect the updated timestamp, and
# update the csv stimulus filenames listed in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650477887014017"
$ws1.Range("B2").Value = "go_stims-16504778869734414.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778869964385.csv"
$ws1.Range("B4").Value = "go_stims-1650477886997439.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778870120165.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778876760848"
$ws2.Range("B2").Value = "TB-16504778875240204.csv"
$ws2.Range("B3").Value = "TB-16504778876550527.csv"
$ws2.Range("B4").Value = "OB-16504778873950171.csv"
$ws2.Range("B5").Value = "OB-16504778872720516.csv"
$ws2.Range("B6").Value = "ZB-match_5-16504778872350192.csv"
$ws2.Range("B7").Value = "ZB-match_3-16504778871070518.csv"
$ws2.Range("B8").Value = "TB-16504778875820217.csv"
$ws2.Range("B9").Value = "ZB-match_5-16504778870340207.csv"
$ws2.Range("B10").Value = "OB-1650477887459022.csv"

# --- Sheet 3: RS_TO --- (name updated, no cell data changes)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778876780167"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778877250175"
$ws4.Range("B2").Value = "MM_stims-16504778876920488.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778876800177.csv"
$ws4.Range("B4").Value = "MM_stims-16504778877080512.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778876930199.csv"
$ws4.Range("B6").Value = "MM_stims-16504778877240498.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778877090197.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778877880495"
$ws5.Range("B2").Value = "SAT_stims-16504778877280247.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778877720494.csv"
$ws5.Range("B4").Value = "SAT_stims-1650477887740017.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778877560515.csv"
